$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45216 -> 45217, i.e. 2023-10-17 -> 2023-10-18) for every data row (2..24).
for ($r = 2; $r -le 24; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45216) {
        $cell.Value2 = 45217
    }
}
